$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2889.4747466532776
$ws.Range("A2").Value = 2794.2232860777544
$ws.Range("A3").Value = 3458.6603638093998
$ws.Range("A4").Value = 1770.7163804517086
$ws.Range("A5").Value = 2870.6052852464877
$ws.Range("A6").Value = 2572.48129521025
$ws.Range("A7").Value = 2848.40187755781
$ws.Range("A8").Value = 2297.917989004287
$ws.Range("A9").Value = 2755.9128469937496
$ws.Range("A10").Value = 2107.97717166034
